# ReportePlanificacion.xlsx - post-planificacion review edit
#
# - "Reporte atrasos" gets a new row (contract 3, delivered late).
# - "Reporte planificación" row 4's dates are corrected and three more
#   planificación rows (2, 3, 4) are added, each with an
#   "Empleados encargados instalación" tag built from combining the
#   installer-id groups "12 ;13" and "4 ;8" (and "4 ;12 ;8 ;13" for the
#   row that used every installer).
# - "Reporte planificación" becomes the active/selected sheet, with M8
#   as the selected cell.

$wb = $excel.ActiveWorkbook

$wsAtrasos = $wb.Worksheets.Item(2)
$wsPlan = $wb.Worksheets.Item(3)

function Set-DateCell($ws, $addr, $serial) {
    $ws.Range($addr).Value = $serial
    $ws.Range($addr).NumberFormat = "yyyy-mm-dd"
}

# ---- "Reporte atrasos" (sheet 2): new row 5 ----
$wsAtrasos.Range("A5").Value = 3
Set-DateCell $wsAtrasos "B5" 42889
Set-DateCell $wsAtrasos "C5" 42892

# ---- "Reporte planificación" (sheet 3) ----

# Row 4: correct the dates / installer-employee counts that were entered
# before the plan was finalised.
$wsPlan.Range("A4").Value = 1
Set-DateCell $wsPlan "B4" 42852
Set-DateCell $wsPlan "C4" 42852
$wsPlan.Range("D4").Value = 1
Set-DateCell $wsPlan "E4" 42853
Set-DateCell $wsPlan "F4" 42853
$wsPlan.Range("G4").Value = 10
Set-DateCell $wsPlan "H4" 42857
Set-DateCell $wsPlan "I4" 42860
$wsPlan.Range("J4").Value = 11
Set-DateCell $wsPlan "K4" 42878
Set-DateCell $wsPlan "L4" 42879
$wsPlan.Range("M4").Value = "4 ;12 ;8 ;13"

# Row 5 (contract 2)
$wsPlan.Range("A5").Value = 2
Set-DateCell $wsPlan "B5" 42852
Set-DateCell $wsPlan "C5" 42852
$wsPlan.Range("D5").Value = 9
Set-DateCell $wsPlan "E5" 42853
Set-DateCell $wsPlan "F5" 42853
$wsPlan.Range("G5").Value = 6
Set-DateCell $wsPlan "H5" 42857
Set-DateCell $wsPlan "I5" 42863
$wsPlan.Range("J5").Value = 7
Set-DateCell $wsPlan "K5" 42880
Set-DateCell $wsPlan "L5" 42885

# Row 6 (contract 3)
$wsPlan.Range("A6").Value = 3
Set-DateCell $wsPlan "B6" 42852
Set-DateCell $wsPlan "C6" 42852
$wsPlan.Range("D6").Value = 9
Set-DateCell $wsPlan "E6" 42853
Set-DateCell $wsPlan "F6" 42853
$wsPlan.Range("G6").Value = 2
Set-DateCell $wsPlan "H6" 42857
Set-DateCell $wsPlan "I6" 42872
$wsPlan.Range("J6").Value = 3
Set-DateCell $wsPlan "K6" 42880
Set-DateCell $wsPlan "L6" 42892
$wsPlan.Range("M6").Value = "12 ;13"

# Row 7 (contract 4)
$wsPlan.Range("A7").Value = 4
Set-DateCell $wsPlan "B7" 42852
Set-DateCell $wsPlan "C7" 42852
$wsPlan.Range("D7").Value = 5
Set-DateCell $wsPlan "E7" 42853
Set-DateCell $wsPlan "F7" 42853
$wsPlan.Range("G7").Value = 2
Set-DateCell $wsPlan "H7" 42863
Set-DateCell $wsPlan "I7" 42870
$wsPlan.Range("J7").Value = 11
Set-DateCell $wsPlan "K7" 42886
Set-DateCell $wsPlan "L7" 42891

# M5 / M7 share the same "4 ;8" tag; set after M6 so the shared-string
# table is built in the same order the workbook ended up with.
$wsPlan.Range("M5").Value = "4 ;8"
$wsPlan.Range("M7").Value = "4 ;8"

# ---- Selection / active sheet: "Reporte planificación" ends up active,
# scrolled/selected at M8 ----
$wsPlan.Activate()
$wsPlan.Range("M8").Select()
